$wb = $excel.ActiveWorkbook

# --- Update the status text from "Ready for handoff" to "In Translation" ---
# This shared string is used on all three sheets: on the "Overview" sheet in
# columns E/F (rows 2-3), and on the per-locale sheets ("zh-cn", "de-de") in
# column C (rows 2-3).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the "Status" columns to fit the shorter text ---
# Overview sheet: columns E and F hold the per-locale status.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de sheets: column C holds the "Status" field.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
